$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.513.68'
$ws.Range('E2').Value = '  -1.53%  '
$ws.Range('D3').Value = '2.439.99'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.95'
$ws.Range('E5').Value = '  -1.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.73'
$ws.Range('E6').Value = '  -2.91%  '
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.532'
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('D9').Value = '2.435.72'
$ws.Range('E9').Value = '  -1.89%  '
$ws.Range('E10').Value = '  -3.62%  '
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.23'
$ws.Range('E12').Value = '  -2.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.353'
$ws.Range('E13').Value = '  -1.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.89'
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000176'
$ws.Range('E15').Value = '  -4.90%  '
$ws.Range('D16').Value = '2.884.22'
$ws.Range('E16').Value = '  +1.13%  '
$ws.Range('D17').Value = '62.378.87'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('D18').Value = '2.433.60'
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.23'
$ws.Range('E19').Value = '  -1.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.25'
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '325.94'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.16'
$ws.Range('E22').Value = '  -1.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.06'
$ws.Range('E23').Value = '  +8.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.13'
$ws.Range('E25').Value = '  -3.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '617.30'
$ws.Range('E26').Value = '  -7.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.92'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').Value = '0.0₃0990'
$ws.Range('E28').Value = '  -5.35%  '
$ws.Range('D29').Value = '2.550.41'
$ws.Range('E29').Value = '  -1.34%  '
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.49'
$ws.Range('E31').Value = '  -1.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.14'
$ws.Range('E32').Value = '  -4.92%  '
$ws.Range('E33').Value = '  -1.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.138'
$ws.Range('E34').Value = '  -4.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.11'
$ws.Range('E35').Value = '  +2.13%  '
$ws.Range('E36').Value = '  -2.90%  '
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.377'
$ws.Range('E38').Value = '  -2.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.70'
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.32'
$ws.Range('E40').Value = '  -3.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '143.98'
$ws.Range('E41').Value = '  -3.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.77'
$ws.Range('E42').Value = '  -5.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.56'
$ws.Range('E43').Value = '  -2.84%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.77'
$ws.Range('E45').Value = '  +0.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '146.91'
$ws.Range('E46').Value = '  -3.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.74'
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.65'
$ws.Range('E48').Value = '  -2.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0527'
$ws.Range('E49').Value = '  -4.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.594'
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0231'
$ws.Range('E51').Value = '  -2.33%  '
